$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 1.62
$ws.Range("I2").Value = 4.75
$ws.Range("U2").Value = 11
$ws.Range("W2").Value = 15
$ws.Range("AF2").Value = 34

# Row 3
$ws.Range("P3").Value = 1.3
$ws.Range("Q3").Value = 3.4
$ws.Range("T3").Value = 9
$ws.Range("U3").Value = 9
$ws.Range("Y3").Value = 21
$ws.Range("AB3").Value = 13
$ws.Range("AE3").Value = 17
$ws.Range("AJ3").Value = 34

# Row 6
$ws.Range("G6").Value = 1.5
$ws.Range("H6").Value = 4.2
$ws.Range("I6").Value = 6
$ws.Range("L6").Value = 1.29
$ws.Range("M6").Value = 3.5
$ws.Range("N6").Value = 1.93
$ws.Range("O6").Value = 1.93
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 1.73
$ws.Range("T6").Value = 6.5
$ws.Range("U6").Value = 6.5
$ws.Range("V6").Value = 8.5
$ws.Range("W6").Value = 10
$ws.Range("Y6").Value = 29
$ws.Range("Z6").Value = 10
$ws.Range("AA6").Value = 8
$ws.Range("AB6").Value = 21
$ws.Range("AC6").Value = 67
$ws.Range("AF6").Value = 29
$ws.Range("AG6").Value = 19
$ws.Range("AH6").Value = 67
$ws.Range("AJ6").Value = 51

# Row 7
$ws.Range("G7").Value = 2.67
$ws.Range("H7").Value = 2.77
$ws.Range("I7").Value = 2.85
$ws.Range("L7").Value = 1.36
$ws.Range("M7").Value = 2.67
$ws.Range("N7").Value = 2.05
$ws.Range("T7").Value = 7.7
$ws.Range("U7").Value = 13.5
$ws.Range("V7").Value = 9.75
$ws.Range("W7").Value = 32
$ws.Range("X7").Value = 24
$ws.Range("Z7").Value = 7.4
$ws.Range("AA7").Value = 5.4
$ws.Range("AB7").Value = 13
$ws.Range("AC7").Value = 65
$ws.Range("AE7").Value = 8.25
$ws.Range("AF7").Value = 14.5
$ws.Range("AG7").Value = 10
$ws.Range("AH7").Value = 37
$ws.Range("AI7").Value = 26

# Row 9
$ws.Range("T9").Value = 16.5
$ws.Range("U9").Value = 35
$ws.Range("Z9").Value = 12.5
$ws.Range("AA9").Value = 7.7
$ws.Range("AB9").Value = 15
$ws.Range("AE9").Value = 8
$ws.Range("AF9").Value = 8.25
$ws.Range("AJ9").Value = 21

# Row 10
$ws.Range("G10").Value = 1.33
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 8
$ws.Range("L10").Value = 1.13
$ws.Range("M10").Value = 6
$ws.Range("T10").Value = 10
$ws.Range("U10").Value = 8.5
$ws.Range("AA10").Value = 10
$ws.Range("AE10").Value = 26
$ws.Range("AG10").Value = 23
$ws.Range("AI10").Value = 51

# Row 11
$ws.Range("I11").Value = 2.8
$ws.Range("K11").Value = 17
$ws.Range("N11").Value = 1.57
$ws.Range("O11").Value = 2.35
$ws.Range("P11").Value = 1.29
$ws.Range("Q11").Value = 3.5
$ws.Range("U11").Value = 15
$ws.Range("Z11").Value = 17
